$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 2.629732666666667
$ws.Range("H2").Value = 7.889198
$ws.Range("I2").Value = 0.07156737804735891
$ws.Range("J2").Value = 0.07156737804735891
$ws.Range("M2").Value = 184.1023456666667
$ws.Range("N2").Value = 552.307037
$ws.Range("O2").Value = 0.9813423747591566
$ws.Range("P2").Value = 0.9813423747591565
$ws.Range("Q2").Value = 484.1399524095918
$ws.Range("R2").Value = 4357.259571686326
$ws.Range("S2").Value = 0.07023210072828152
$ws.Range("T2").Value = 0.07023210072828152

$ws.Range("G3").Value = 2.629732666666667
$ws.Range("H3").Value = 7.889198
$ws.Range("I3").Value = 0.07156737804735891
$ws.Range("J3").Value = 0.07156737804735891
$ws.Range("O3").Value = 0.002303378255889225
$ws.Range("P3").Value = 0.002303378255889224
$ws.Range("Q3").Value = 1.136359203342445
$ws.Range("R3").Value = 10.227232830082
$ws.Range("S3").Value = 0.0001648467424252904
$ws.Range("T3").Value = 0.0001648467424252903

$ws.Range("G4").Value = 2.629732666666667
$ws.Range("H4").Value = 7.889198
$ws.Range("I4").Value = 0.07156737804735891
$ws.Range("J4").Value = 0.07156737804735891
$ws.Range("M4").Value = 1.367901
$ws.Range("N4").Value = 4.103703
$ws.Range("O4").Value = 0.007291483500193526
$ws.Range("P4").Value = 0.007291483500193526
$ws.Range("Q4").Value = 3.597213944466001
$ws.Range("R4").Value = 32.374925500194
$ws.Range("S4").Value = 0.0005218323561844299
$ws.Range("T4").Value = 0.0005218323561844299

$ws.Range("G5").Value = 2.629732666666667
$ws.Range("H5").Value = 7.889198
$ws.Range("I5").Value = 0.07156737804735891
$ws.Range("J5").Value = 0.07156737804735891
$ws.Range("M5").Value = 1.700197666666667
$ws.Range("N5").Value = 5.100593
$ws.Range("O5").Value = 0.009062763484760617
$ws.Range("P5").Value = 0.009062763484760615
$ws.Range("Q5").Value = 4.471065343823779
$ws.Range("R5").Value = 40.239588094414
$ws.Range("S5").Value = 0.0006485982204676629
$ws.Range("T5").Value = 0.0006485982204676628

$ws.Range("I6").Value = 0.493312042610523
$ws.Range("J6").Value = 0.493312042610523
$ws.Range("M6").Value = 184.1023456666667
$ws.Range("N6").Value = 552.307037
$ws.Range("O6").Value = 0.9813423747591566
$ws.Range("P6").Value = 0.9813423747591565
$ws.Range("Q6").Value = 3337.163877576914
$ws.Range("R6").Value = 30034.47489819222
$ws.Range("S6").Value = 0.4841080113927009
$ws.Range("T6").Value = 0.4841080113927009

$ws.Range("I7").Value = 0.493312042610523
$ws.Range("J7").Value = 0.493312042610523
$ws.Range("O7").Value = 0.002303378255889225
$ws.Range("P7").Value = 0.002303378255889224
$ws.Range("S7").Value = 0.001136284232317378
$ws.Range("T7").Value = 0.001136284232317377

$ws.Range("I8").Value = 0.493312042610523
$ws.Range("J8").Value = 0.493312042610523
$ws.Range("M8").Value = 1.367901
$ws.Range("N8").Value = 4.103703
$ws.Range("O8").Value = 0.007291483500193526
$ws.Range("P8").Value = 0.007291483500193526
$ws.Range("Q8").Value = 24.795500506911
$ws.Range("R8").Value = 223.159504562199
$ws.Range("S8").Value = 0.003596976619141394
$ws.Range("T8").Value = 0.003596976619141394

$ws.Range("I9").Value = 0.493312042610523
$ws.Range("J9").Value = 0.493312042610523
$ws.Range("M9").Value = 1.700197666666667
$ws.Range("N9").Value = 5.100593
$ws.Range("O9").Value = 0.009062763484760617
$ws.Range("P9").Value = 0.009062763484760615
$ws.Range("Q9").Value = 30.81893507328545
$ws.Range("R9").Value = 277.370415659569
$ws.Range("S9").Value = 0.004470770366363321
$ws.Range("T9").Value = 0.004470770366363321

$ws.Range("G10").Value = 7.550656333333333
$ws.Range("H10").Value = 22.651969
$ws.Range("I10").Value = 0.2054888252189962
$ws.Range("J10").Value = 0.2054888252189962
$ws.Range("M10").Value = 184.1023456666667
$ws.Range("N10").Value = 552.307037
$ws.Range("O10").Value = 0.9813423747591566
$ws.Range("P10").Value = 0.9813423747591565
$ws.Range("Q10").Value = 1390.093542289539
$ws.Range("R10").Value = 12510.84188060585
$ws.Range("S10").Value = 0.201654891726879
$ws.Range("T10").Value = 0.201654891726879

$ws.Range("G11").Value = 7.550656333333333
$ws.Range("H11").Value = 22.651969
$ws.Range("I11").Value = 0.2054888252189962
$ws.Range("J11").Value = 0.2054888252189962
$ws.Range("O11").Value = 0.002303378255889225
$ws.Range("P11").Value = 0.002303378255889224
$ws.Range("Q11").Value = 3.262787097874556
$ws.Range("R11").Value = 29.365083880871
$ws.Range("S11").Value = 0.0004733184918376572
$ws.Range("T11").Value = 0.0004733184918376572

$ws.Range("G12").Value = 7.550656333333333
$ws.Range("H12").Value = 22.651969
$ws.Range("I12").Value = 0.2054888252189962
$ws.Range("J12").Value = 0.2054888252189962
$ws.Range("M12").Value = 1.367901
$ws.Range("N12").Value = 4.103703
$ws.Range("O12").Value = 0.007291483500193526
$ws.Range("P12").Value = 0.007291483500193526
$ws.Range("Q12").Value = 10.328550349023
$ws.Range("R12").Value = 92.95695314120701
$ws.Range("S12").Value = 0.001498318378558462
$ws.Range("T12").Value = 0.001498318378558462

$ws.Range("G13").Value = 7.550656333333333
$ws.Range("H13").Value = 22.651969
$ws.Range("I13").Value = 0.2054888252189962
$ws.Range("J13").Value = 0.2054888252189962
$ws.Range("M13").Value = 1.700197666666667
$ws.Range("N13").Value = 5.100593
$ws.Range("O13").Value = 0.009062763484760617
$ws.Range("P13").Value = 0.009062763484760615
$ws.Range("Q13").Value = 12.83760827973522
$ws.Range("R13").Value = 115.538474517617
$ws.Range("S13").Value = 0.001862296621721075
$ws.Range("T13").Value = 0.001862296621721075

$ws.Range("G14").Value = 8.437784666666667
$ws.Range("H14").Value = 25.313354
$ws.Range("I14").Value = 0.2296317541231219
$ws.Range("J14").Value = 0.2296317541231219
$ws.Range("M14").Value = 184.1023456666667
$ws.Range("N14").Value = 552.307037
$ws.Range("O14").Value = 0.9813423747591566
$ws.Range("P14").Value = 0.9813423747591565
$ws.Range("Q14").Value = 1553.415949363567
$ws.Range("R14").Value = 13980.7435442721
$ws.Range("S14").Value = 0.2253473709112952
$ws.Range("T14").Value = 0.2253473709112951

$ws.Range("G15").Value = 8.437784666666667
$ws.Range("H15").Value = 25.313354
$ws.Range("I15").Value = 0.2296317541231219
$ws.Range("J15").Value = 0.2296317541231219
$ws.Range("O15").Value = 0.002303378255889225
$ws.Range("P15").Value = 0.002303378255889224
$ws.Range("Q15").Value = 3.646132697565112
$ws.Range("R15").Value = 32.815194278086
$ws.Range("S15").Value = 0.0005289287893088997
$ws.Range("T15").Value = 0.0005289287893088996

$ws.Range("G16").Value = 8.437784666666667
$ws.Range("H16").Value = 25.313354
$ws.Range("I16").Value = 0.2296317541231219
$ws.Range("J16").Value = 0.2296317541231219
$ws.Range("M16").Value = 1.367901
$ws.Range("N16").Value = 4.103703
$ws.Range("O16").Value = 0.007291483500193526
$ws.Range("P16").Value = 0.007291483500193526
$ws.Range("Q16").Value = 11.542054083318
$ws.Range("R16").Value = 103.878486749862
$ws.Range("S16").Value = 0.00167435614630924
$ws.Range("T16").Value = 0.00167435614630924

$ws.Range("G17").Value = 8.437784666666667
$ws.Range("H17").Value = 25.313354
$ws.Range("I17").Value = 0.2296317541231219
$ws.Range("J17").Value = 0.2296317541231219
$ws.Range("M17").Value = 1.700197666666667
$ws.Range("N17").Value = 5.100593
$ws.Range("O17").Value = 0.009062763484760617
$ws.Range("P17").Value = 0.009062763484760615
$ws.Range("Q17").Value = 14.34590180210245
$ws.Range("R17").Value = 129.113116218922
$ws.Range("S17").Value = 0.002081098276208557
$ws.Range("T17").Value = 0.002081098276208557
